$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all Timestamp values in column A (rows 2-97) forward by 8 days
for ($r = 2; $r -le 97; $r++) {
    $cur = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 1).Value = $cur + 8
}

# Update Notified Production (MW) values in column B for the affected rows
$newB = @{
    18 = 2
    19 = 2
    20 = 3
    21 = 4
    22 = 41
    23 = 46
    24 = 54
    25 = 65
    26 = 207
    27 = 223
    28 = 241
    29 = 262
    30 = 550
    31 = 577
    32 = 605
    33 = 635
    34 = 888
    35 = 917
    36 = 948
    37 = 978
    38 = 1201
    39 = 1229
    40 = 1256
    41 = 1283
    42 = 1460
    43 = 1480
    44 = 1497
    45 = 1512
    46 = 1584
    47 = 1591
    48 = 1600
    49 = 1609
    50 = 1634
    51 = 1638
    52 = 1644
    53 = 1649
    54 = 1643
    55 = 1643
    56 = 1640
    57 = 1634
    58 = 1542
    59 = 1528
    60 = 1511
    61 = 1493
    62 = 1311
    63 = 1285
    64 = 1258
    65 = 1229
    66 = 1017
    67 = 986
    68 = 955
    69 = 924
    70 = 595
    71 = 564
    72 = 534
    73 = 504
    74 = 232
    75 = 209
    76 = 189
    77 = 171
    78 = 54
    79 = 44
    80 = 37
    81 = 33
    82 = 24
    83 = 24
    84 = 24
    85 = 24
    86 = 9
    87 = 9
    88 = 9
    89 = 9
}
foreach ($r in $newB.Keys) {
    $ws.Cells.Item($r, 2).Value = $newB[$r]
}

